# Applies the scheduled text replacements for the addition/subtraction worksheet.
$d = $word.ActiveDocument

function Replace-ExactText($doc, $old, $new) {
    # MatchCase=$true, MatchWholeWord=$true, Forward=$true, Wrap=wdFindContinue(1),
    # Format=$false, Replace=wdReplaceOne(1)
    $found = $doc.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 1)
    if (-not $found) {
        throw "Could not find expected text: $old"
    }
}

Replace-ExactText $d "2024-12-30 Monday" "2024-12-31 Tuesday"
Replace-ExactText $d "1+79=80" "91-70=21"
Replace-ExactText $d "90-24=66" "2+47=49"
Replace-ExactText $d "14+46=60" "59+3=62"
Replace-ExactText $d "61-49=12" "19+49=68"
Replace-ExactText $d "9-6=3" "80-2=78"
Replace-ExactText $d "17+52=69" "67+27=94"
Replace-ExactText $d "90-89=1" "94-71=23"
Replace-ExactText $d "60-36=24" "11+43=54"
Replace-ExactText $d "82-3=79" "98-42=56"
Replace-ExactText $d "78-58=20" "6+3=9"
Replace-ExactText $d "90+9=99" "92-26=66"
Replace-ExactText $d "34+18=52" "26+43=69"
Replace-ExactText $d "97-72=25" "37+24=61"
Replace-ExactText $d "47-38=9" "44-30=14"
Replace-ExactText $d "70-10=60" "27+5=32"
Replace-ExactText $d "96-46=50" "38-28=10"
Replace-ExactText $d "24+41=65" "55+8=63"
Replace-ExactText $d "73-44=29" "58-50=8"
Replace-ExactText $d "88-23=65" "93-50=43"
Replace-ExactText $d "52+37=89" "39+8=47"
Replace-ExactText $d "3+71=74" "24+13=37"
Replace-ExactText $d "96-5=91" "45+45=90"
Replace-ExactText $d "45+41=86" "8+0=8"
Replace-ExactText $d "1+55=56" "6+45=51"
Replace-ExactText $d "42-20=22" "70+29=99"
Replace-ExactText $d "57-18=39" "49-40=9"
Replace-ExactText $d "87-11=76" "59-33=26"
Replace-ExactText $d "79+9=88" "11-4=7"
Replace-ExactText $d "44+10=54" "5+86=91"
Replace-ExactText $d "46-34=12" "30+5=35"
Replace-ExactText $d "64-61=3" "38-24=14"
Replace-ExactText $d "93-33=60" "59+27=86"
Replace-ExactText $d "42-14=28" "26-2=24"
Replace-ExactText $d "55+31=86" "49+37=86"
Replace-ExactText $d "83-82=1" "82-6=76"
Replace-ExactText $d "82-65=17" "65+27=92"
Replace-ExactText $d "76+10=86" "5+11=16"
Replace-ExactText $d "93-0=93" "58+29=87"
Replace-ExactText $d "85-39=46" "98-42=56"
Replace-ExactText $d "77-16=61" "24+69=93"
Replace-ExactText $d "16+58=74" "90-67=23"
Replace-ExactText $d "52+3=55" "6-2=4"
Replace-ExactText $d "54-1=53" "82+5=87"
Replace-ExactText $d "81-31=50" "49-31=18"
Replace-ExactText $d "81-19=62" "75-54=21"
Replace-ExactText $d "51+44=95" "63-40=23"
Replace-ExactText $d "55+44=99" "86-34=52"
Replace-ExactText $d "36+6=42" "57+37=94"
Replace-ExactText $d "79-7=72" "56-35=21"
Replace-ExactText $d "29-26=3" "49-40=9"
Replace-ExactText $d "41-2=39" "59-1=58"
Replace-ExactText $d "55+35=90" "29+26=55"
Replace-ExactText $d "51-11=40" "56-2=54"
Replace-ExactText $d "1+16=17" "81-14=67"
Replace-ExactText $d "69-10=59" "94-55=39"
Replace-ExactText $d "97-85=12" "12+37=49"
Replace-ExactText $d "42+5=47" "31-18=13"
Replace-ExactText $d "37+21=58" "46-37=9"
Replace-ExactText $d "18+80=98" "37-32=5"
Replace-ExactText $d "62-22=40" "72-69=3"
Replace-ExactText $d "36+14=50" "26-13=13"
Replace-ExactText $d "41+47=88" "56+43=99"
Replace-ExactText $d "84+0=84" "28+27=55"
Replace-ExactText $d "17+34=51" "16+29=45"
Replace-ExactText $d "12+87=99" "85-23=62"
Replace-ExactText $d "74-30=44" "23+19=42"
Replace-ExactText $d "66-29=37" "41+22=63"
Replace-ExactText $d "39+15=54" "72-49=23"
Replace-ExactText $d "79+13=92" "44-5=39"
Replace-ExactText $d "86-18=68" "37+13=50"
Replace-ExactText $d "33-5=28" "29+69=98"
Replace-ExactText $d "24-2=22" "76+8=84"
Replace-ExactText $d "74+23=97" "9+77=86"
Replace-ExactText $d "44+38=82" "2+21=23"
Replace-ExactText $d "40-10=30" "72-42=30"
Replace-ExactText $d "92-5=87" "72+4=76"
Replace-ExactText $d "19+31=50" "44-13=31"
Replace-ExactText $d "34-7=27" "50-15=35"
Replace-ExactText $d "64-36=28" "16+71=87"
Replace-ExactText $d "62+34=96" "91-59=32"
Replace-ExactText $d "31+68=99" "8+74=82"
Replace-ExactText $d "50-37=13" "49+11=60"
Replace-ExactText $d "96-25=71" "42-36=6"
Replace-ExactText $d "80+9=89" "90-58=32"
Replace-ExactText $d "97-25=72" "94-67=27"
Replace-ExactText $d "86-68=18" "76-32=44"
Replace-ExactText $d "44-21=23" "40+51=91"
Replace-ExactText $d "43-25=18" "3+65=68"
Replace-ExactText $d "65+20=85" "8+88=96"
Replace-ExactText $d "37+19=56" "35+62=97"
Replace-ExactText $d "44-2=42" "97-17=80"
Replace-ExactText $d "18-0=18" "11+79=90"
Replace-ExactText $d "78-57=21" "73-60=13"
Replace-ExactText $d "96-26=70" "59-4=55"
Replace-ExactText $d "51+17=68" "84+2=86"
Replace-ExactText $d "22+62=84" "8+38=46"
Replace-ExactText $d "48-45=3" "90-55=35"
Replace-ExactText $d "62-54=8" "34-22=12"
Replace-ExactText $d "64-9=55" "69-31=38"
Replace-ExactText $d "36-5=31" "59-42=17"
